$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "1.00", "0.482") are preserved verbatim instead of being
# coerced into numbers by Excel's COM value-assignment heuristics.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.643.72'
$ws.Range("E2").Value = '  -3.98%  '
$ws.Range("D3").Value = '3.455.90'
$ws.Range("E3").Value = '  -4.15%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '601.96'
$ws.Range("E5").Value = '  -4.22%  '
$ws.Range("D6").Value = '147.65'
$ws.Range("E6").Value = '  -6.83%  '
$ws.Range("D7").Value = '3.450.94'
$ws.Range("E7").Value = '  -4.28%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").Value = '  -2.70%  '
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  -4.92%  '
$ws.Range("D11").Value = '7.45'
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("E12").Value = '  -4.03%  '
$ws.Range("D13").Value = '0.0000212'
$ws.Range("E13").Value = '  -7.57%  '
$ws.Range("D14").Value = '31.49'
$ws.Range("E14").Value = '  -6.06%  '
$ws.Range("D15").Value = '4.038.26'
$ws.Range("E15").Value = '  -4.45%  '
$ws.Range("D16").Value = '3.449.65'
$ws.Range("E16").Value = '  -4.22%  '
$ws.Range("D17").Value = '66.767.39'
$ws.Range("E17").Value = '  -4.17%  '
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  -4.51%  '
$ws.Range("D20").Value = '15.23'
$ws.Range("E20").Value = '  -5.59%  '
$ws.Range("D21").Value = '9.99'
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").Value = '438.71'
$ws.Range("E22").Value = '  -5.15%  '
$ws.Range("E23").Value = '  -5.46%  '
$ws.Range("D24").Value = '78.10'
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '3.594.66'
$ws.Range("E26").Value = '  -4.36%  '
$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  -11.22%  '
$ws.Range("D28").Value = '9.80'
$ws.Range("E28").Value = '  -8.23%  '
$ws.Range("D29").Value = '8.37'
$ws.Range("E29").Value = '  -9.73%  '
$ws.Range("D30").Value = '2.47'
$ws.Range("E30").Value = '  -6.25%  '
$ws.Range("D31").Value = '1.60'
$ws.Range("E31").Value = '  -7.06%  '
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").Value = '0.165'
$ws.Range("E33").Value = '  -6.17%  '
$ws.Range("D34").Value = '25.37'
$ws.Range("E34").Value = '  -4.30%  '
$ws.Range("D35").Value = '6.09'
$ws.Range("E35").Value = '  -6.89%  '
$ws.Range("D36").Value = '3.447.88'
$ws.Range("E36").Value = '  -4.54%  '
$ws.Range("D37").Value = '1.80'
$ws.Range("E37").Value = '  -7.73%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '7.89'
$ws.Range("E39").Value = '  -7.07%  '
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").Value = '173.50'
$ws.Range("E41").Value = '  -3.61%  '
$ws.Range("D43").Value = '0.0885'
$ws.Range("E43").Value = '  -4.25%  '
$ws.Range("D44").Value = '5.37'
$ws.Range("E44").Value = '  -5.19%  '
$ws.Range("D45").Value = '0.880'
$ws.Range("E45").Value = '  -3.66%  '
$ws.Range("D46").Value = '29.06'
$ws.Range("E46").Value = '  -9.04%  '
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("D48").Value = '1.23'
$ws.Range("E48").Value = '  -10.95%  '
$ws.Range("D49").Value = '7.47'
$ws.Range("E49").Value = '  -4.37%  '
$ws.Range("D50").Value = '2.45'
$ws.Range("E50").Value = '  -10.93%  '
$ws.Range("D51").Value = '0.987'
$ws.Range("E51").Value = '  -5.08%  '

# Restore the default "Normal" style on column D so we don't leave
# a stray text-format override on cells that didn't have one before.
$ws.Range("D2:D51").Style = "Normal"
